$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value2 = 89124839
$ws.Range('B2').Value2 = 95525
$ws.Range('D2').Value2 = 'LC'
$ws.Range('E2').Value2 = 221941
$ws.Range('F2').Value2 = 'Plattlummer'
$ws.Range('G2').Value2 = 'Lycopodium complanatum'
$ws.Range('H2').Value2 = 'L.'
$ws.Range('Q2').Value2 = 533753.1171408413
$ws.Range('R2').Value2 = 6903109.937925656

# Row 3
$ws.Range('A3').Value2 = 89124842
$ws.Range('B3').Value2 = 95519
$ws.Range('D3').Value2 = 'LC'
$ws.Range('E3').Value2 = 221945
$ws.Range('F3').Value2 = 'Revlummer'
$ws.Range('G3').Value2 = 'Lycopodium annotinum'
$ws.Range('H3').Value2 = 'L.'
$ws.Range('Q3').Value2 = 533669.7937915208
$ws.Range('R3').Value2 = 6903090.933150688

# Row 4
$ws.Range('A4').Value2 = 89124843
$ws.Range('B4').Value2 = 77506
$ws.Range('D4').Value2 = 'NT'
$ws.Range('E4').Value2 = 6425
$ws.Range('F4').Value2 = 'Garnlav'
$ws.Range('G4').Value2 = 'Alectoria sarmentosa'
$ws.Range('H4').Value2 = '(Ach.) Ach.'
$ws.Range('Q4').Value2 = 533659.0313536879
$ws.Range('R4').Value2 = 6903094.08653159

# Row 5
$ws.Range('A5').Value2 = 89124841
$ws.Range('B5').Value2 = 96237
$ws.Range('D5').Value2 = 'LC'
$ws.Range('E5').Value2 = 220093
$ws.Range('F5').Value2 = 'Korallrot'
$ws.Range('G5').Value2 = 'Corallorhiza trifida'
$ws.Range('H5').Value2 = 'Châtel.'
$ws.Range('Q5').Value2 = 533677.0621316924
$ws.Range('R5').Value2 = 6903063.984857111

# Row 6
$ws.Range('A6').Value2 = 89124858
$ws.Range('B6').Value2 = 77259
$ws.Range('D6').Value2 = 'NT'
$ws.Range('E6').Value2 = 228912
$ws.Range('F6').Value2 = 'Mörk kolflarnlav'
$ws.Range('G6').Value2 = 'Carbonicola myrmecina'
$ws.Range('H6').Value2 = '(Ach.) Bendiksby & Timdal'
$ws.Range('Q6').Value2 = 534185.1622770416
$ws.Range('R6').Value2 = 6902557.999860712

# Row 7
$ws.Range('A7').Value2 = 89124851
$ws.Range('B7').Value2 = 77258
$ws.Range('D7').Value2 = 'NT'
$ws.Range('E7').Value2 = 6446
$ws.Range('F7').Value2 = 'Kolflarnlav'
$ws.Range('G7').Value2 = 'Carbonicola anthracophila'
$ws.Range('H7').Value2 = '(Nyl.) Bendiksby & Timdal'
$ws.Range('Q7').Value2 = 533999.1884141648
$ws.Range('R7').Value2 = 6902813.777029264

# Row 8
$ws.Range('A8').Value2 = 89124856
$ws.Range('B8').Value2 = 77506
$ws.Range('D8').Value2 = 'NT'
$ws.Range('E8').Value2 = 6425
$ws.Range('F8').Value2 = 'Garnlav'
$ws.Range('G8').Value2 = 'Alectoria sarmentosa'
$ws.Range('H8').Value2 = '(Ach.) Ach.'
$ws.Range('Q8').Value2 = 534125.9184093268
$ws.Range('R8').Value2 = 6902602.129391655

# Row 9
$ws.Range('A9').Value2 = 89124850
$ws.Range('B9').Value2 = 77506
$ws.Range('D9').Value2 = 'NT'
$ws.Range('E9').Value2 = 6425
$ws.Range('F9').Value2 = 'Garnlav'
$ws.Range('G9').Value2 = 'Alectoria sarmentosa'
$ws.Range('H9').Value2 = '(Ach.) Ach.'
$ws.Range('Q9').Value2 = 533977.1051462417
$ws.Range('R9').Value2 = 6902828.929063975

# Row 10
$ws.Range('A10').Value2 = 89124852
$ws.Range('B10').Value2 = 77506
$ws.Range('D10').Value2 = 'NT'
$ws.Range('E10').Value2 = 6425
$ws.Range('F10').Value2 = 'Garnlav'
$ws.Range('G10').Value2 = 'Alectoria sarmentosa'
$ws.Range('H10').Value2 = '(Ach.) Ach.'
$ws.Range('Q10').Value2 = 534000.1686273545
$ws.Range('R10').Value2 = 6902809.127754148

# Row 11
$ws.Range('A11').Value2 = 89124855
$ws.Range('B11').Value2 = 76909
$ws.Range('D11').Value2 = 'NT'
$ws.Range('E11').Value2 = 6437
$ws.Range('F11').Value2 = 'Blanksvart spiklav'
$ws.Range('G11').Value2 = 'Calicium denigratum'
$ws.Range('H11').Value2 = '(Vain.) Tibell'
$ws.Range('Q11').Value2 = 534081.1437023395
$ws.Range('R11').Value2 = 6902692.067965358

# Row 12
$ws.Range('A12').Value2 = 89124859
$ws.Range('B12').Value2 = 77506
$ws.Range('D12').Value2 = 'NT'
$ws.Range('E12').Value2 = 6425
$ws.Range('F12').Value2 = 'Garnlav'
$ws.Range('G12').Value2 = 'Alectoria sarmentosa'
$ws.Range('H12').Value2 = '(Ach.) Ach.'
$ws.Range('Q12').Value2 = 534178.1517579975
$ws.Range('R12').Value2 = 6902604.989653899
